# Fixed free camera issue by regenerating matrix from euler angles each time.
# This script reproduces the corresponding spreadsheet edit:
#  - Removes the "Add camera controls to model viewer" / "Model Viewer" task row
#    (and its associated comment about the camera rolling around).
#  - Moves the "On-screen debug text" task down the list (after "Compiled
#    shaders") and bumps its estimate from 14 to 21.
#  - Shifts everything else up to fill the gap, keeping the existing blank
#    separator row before the second table ("Scene Exporter" section).
#  - Re-anchors the remaining three cell comments to their new rows.
#  - Updates the active selection from D12 to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture comment text (in original cell order) before we touch data,
#    so we can recreate them at their new locations afterwards.
# ---------------------------------------------------------------------
$commentB5 = $ws.Range("B5").Comment.Text()
$commentB13 = $ws.Range("B13").Comment.Text()
$commentB18 = $ws.Range("B18").Comment.Text()

# Remove all four existing comments (B2, B5, B13, B18). B2's comment is not
# recreated -- it belonged to the task row being deleted entirely.
$ws.Range("B2").Comment.Delete()
$ws.Range("B5").Comment.Delete()
$ws.Range("B13").Comment.Delete()
$ws.Range("B18").Comment.Delete()

# ---------------------------------------------------------------------
# 2. Clear out the old table body (rows 2-22) and rewrite it with the new
#    row order / values.
# ---------------------------------------------------------------------
$ws.Range("A2:C22").ClearContents()

$rows = @(
    @("Engine",         "Point lights", 21),
    @("Engine",         "Point lights in model chain", 7),
    @("Engine",         "Bounding volume occlusion", 7),
    @("Rorn Pool",      "Build basic app", 3),
    @("Rorn Pool",      "Game initialisation (position of balls, etc)", 14),
    @("Engine",         "Basics of a physics engine", 21),
    @("Rorn Pool",      "Player can move the cue ball (when appropriate)", 5),
    @("Rorn Pool",      "Player can move cue", 5),
    @("Rorn Pool",      "Player can take a shot", 10),
    @("Engine",         "Compiled shaders", 5),
    @("Engine",         "On-screen debug text", 21)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Row 13 stays blank (separator), just like the old row 14 did.

$rows2 = @(
    @("Scene Exporter",  "Make exporter a GUP.  Build UI & hook data into the max files", 14),
    @("Scene Exporter",  "Make the path from Max->Model Viewer seamless", 4),
    @("Model Compiler",  "Model compiler refactoring - wait till we have two surface formats implemented", 21),
    @("Model Compiler",  "Error handling strategy in Model Compiler - possibly use xsd?", 21),
    @("Maths",           "Complete the Rorn Maths library", 35),
    @("Engine",          "Revise, understand and document the view and projection matrix builds", 7),
    @("Engine",          "Add full screen support", 3),
    @("Scene Exporter",  "Get 64-bit scene exporter working", 7)
)

$r = 14
foreach ($row in $rows2) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Recreate the surviving comments at their new (shifted) locations.
# ---------------------------------------------------------------------
$ws.Range("B3").AddComment($commentB5)
$ws.Range("B11").AddComment($commentB13)
$ws.Range("B17").AddComment($commentB18)

# ---------------------------------------------------------------------
# 4. Update the selected / active cell to match the new state.
# ---------------------------------------------------------------------
$ws.Range("C12").Select()
